$d = $word.ActiveDocument

# Find the paragraph that holds "Lo-lo-lo" and insert a brand-new
# paragraph ("1111") directly after it, matching the commit's intent.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Lo-lo-lo") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $r = $d.Paragraphs.Item($targetIndex).Range
    $r.Collapse(0)             # wdCollapseEnd
    $r.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "1111"
    Write-Output "Inserted new paragraph '1111' after paragraph $targetIndex."
} else {
    Write-Output "Could not locate the 'Lo-lo-lo' paragraph; no changes made."
}
